# Slide 42 ("Code Generation for CompoundStmt"): in the grammar-rule
# paragraph "compoundStmt = "{" statements "}" .", the runs
#   ' = ', '"{" statements ', '"}" .'
# are collapsed into a single run (keeping the formatting of the
# ' = ' run) with combined text ' = "{" statements "}" .'.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(42)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$full = $tr.Text
$oldSuffix = ' = "{" statements "}" .'
$idx0 = $full.IndexOf($oldSuffix)

if ($idx0 -ge 0) {
    $start1 = $idx0 + 1
    $target = $tr.Characters($start1, $oldSuffix.Length)
    $target.Text = $oldSuffix
}
